# Auto update: 2025-12-06 21:20:02
# Refresh the analysis table (rows 2-7, columns D-O) on the active sheet with the
# latest pull of market / model data. Column A (date), B (name/ticker-info) and
# C (ticker) keep referencing the same underlying values; only the numeric /
# categorical figures for this run change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : SamsungElec / 005930.KS
$ws.Range("B2").Value = "SamsungElec"
$ws.Range("C2").Value = "005930.KS"
$ws.Range("D2").Value = 108400
$ws.Range("E2").Value = 61.6
$ws.Range("F2").Value = 7.86
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 76
$ws.Range("K2").Value = 54.7
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 52.28493729186943
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3 : 058470.KS
$ws.Range("B3").Value = "058470.KS,0P0000ASU1,98886"
$ws.Range("C3").Value = "058470.KS"
$ws.Range("D3").Value = 65400
$ws.Range("E3").Value = 61.7
$ws.Range("F3").Value = -4.25
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 51.9
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 52.28493729186943
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 4 : 403870.KS
$ws.Range("D4").Value = 30300
$ws.Range("E4").Value = 43.3
$ws.Range("F4").Value = 0.17
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 53
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 50.7
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 52.28493729186943
$ws.Range("O4").Value = "⚪ 중립 구간"

# Row 5 : SK hynix / 000660.KS (only K/M/N/O change this run)
$ws.Range("K5").Value = 48.1
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 52.28493729186943
$ws.Range("O5").Value = "⚪ 중립 구간"

# Row 6 : 240810.KS
$ws.Range("B6").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C6").Value = "240810.KS"
$ws.Range("D6").Value = 61800
$ws.Range("E6").Value = 38
$ws.Range("F6").Value = 0.82
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 53
$ws.Range("I6").Value = 53
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 42.9
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 52.28493729186943
$ws.Range("O6").Value = "⚪ 중립 구간"

# Row 7 : DB HiTek / 000990.KS
$ws.Range("B7").Value = "DB HiTek"
$ws.Range("C7").Value = "000990.KS"
$ws.Range("D7").Value = 64800
$ws.Range("E7").Value = 33.9
$ws.Range("F7").Value = 1.89
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = 46
$ws.Range("J7").Value = 63
$ws.Range("K7").Value = 40.1
$ws.Range("M7").Value = "⛔ 관망하십시오."
$ws.Range("N7").Value = 52.28493729186943
$ws.Range("O7").Value = "⚪ 중립 구간"
